$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update correct/total marks in the marksheet
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/140"
